$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 5.021024000000001
$ws.Range("H2").Value = 15.063072
$ws.Range("I2").Value = 0.805437917812573
$ws.Range("J2").Value = 0.805437917812573
$ws.Range("M2").Value = 0.5550926666666666
$ws.Range("N2").Value = 1.665278
$ws.Range("O2").Value = 0.1208967663154349
$ws.Range("P2").Value = 0.1208967663154349
$ws.Range("Q2").Value = 2.787133601557333
$ws.Range("R2").Value = 25.084202414016
$ws.Range("S2").Value = 0.09737483973137707
$ws.Range("T2").Value = 0.09737483973137707
$ws.Range("G3").Value = 5.021024000000001
$ws.Range("H3").Value = 15.063072
$ws.Range("I3").Value = 0.805437917812573
$ws.Range("J3").Value = 0.805437917812573
$ws.Range("O3").Value = 0.7377399926530269
$ws.Range("P3").Value = 0.7377399926530268
$ws.Range("Q3").Value = 17.00773300562133
$ws.Range("R3").Value = 153.069597050592
$ws.Range("S3").Value = 0.5942037635695169
$ws.Range("T3").Value = 0.5942037635695169
$ws.Range("G4").Value = 5.021024000000001
$ws.Range("H4").Value = 15.063072
$ws.Range("I4").Value = 0.805437917812573
$ws.Range("J4").Value = 0.805437917812573
$ws.Range("M4").Value = 0.5311786666666667
$ws.Range("N4").Value = 1.593536
$ws.Range("O4").Value = 0.1156884012202364
$ws.Range("P4").Value = 0.1156884012202364
$ws.Range("Q4").Value = 2.667060833621334
$ws.Range("R4").Value = 24.003547502592
$ws.Range("S4").Value = 0.09317982499389275
$ws.Range("T4").Value = 0.09317982499389274
$ws.Range("G5").Value = 5.021024000000001
$ws.Range("H5").Value = 15.063072
$ws.Range("I5").Value = 0.805437917812573
$ws.Range("J5").Value = 0.805437917812573
$ws.Range("M5").Value = 0.117885
$ws.Range("N5").Value = 0.353655
$ws.Range("O5").Value = 0.02567483981130185
$ws.Range("P5").Value = 0.02567483981130185
$ws.Range("Q5").Value = 0.59190341424
$ws.Range("R5").Value = 5.32713072816
$ws.Range("S5").Value = 0.02067948951778632
$ws.Range("T5").Value = 0.02067948951778632
$ws.Range("I6").Value = 0.0482448215850983
$ws.Range("J6").Value = 0.0482448215850983
$ws.Range("M6").Value = 0.5550926666666666
$ws.Range("N6").Value = 1.665278
$ws.Range("O6").Value = 0.1208967663154349
$ws.Range("P6").Value = 0.1208967663154349
$ws.Range("Q6").Value = 0.1669461548397777
$ws.Range("R6").Value = 1.502515393558
$ws.Range("S6").Value = 0.005832642921103477
$ws.Range("T6").Value = 0.005832642921103477
$ws.Range("I7").Value = 0.0482448215850983
$ws.Range("J7").Value = 0.0482448215850983
$ws.Range("O7").Value = 0.7377399926530269
$ws.Range("P7").Value = 0.7377399926530268
$ws.Range("S7").Value = 0.03559213432173701
$ws.Range("T7").Value = 0.03559213432173701
$ws.Range("I8").Value = 0.0482448215850983
$ws.Range("J8").Value = 0.0482448215850983
$ws.Range("M8").Value = 0.5311786666666667
$ws.Range("N8").Value = 1.593536
$ws.Range("O8").Value = 0.1156884012202364
$ws.Range("P8").Value = 0.1156884012202364
$ws.Range("Q8").Value = 0.1597539316551111
$ws.Range("R8").Value = 1.437785384896
$ws.Range("S8").Value = 0.005581366276335574
$ws.Range("T8").Value = 0.005581366276335574
$ws.Range("I9").Value = 0.0482448215850983
$ws.Range("J9").Value = 0.0482448215850983
$ws.Range("M9").Value = 0.117885
$ws.Range("N9").Value = 0.353655
$ws.Range("O9").Value = 0.02567483981130185
$ws.Range("P9").Value = 0.02567483981130185
$ws.Range("Q9").Value = 0.035454345995
$ws.Range("R9").Value = 0.319089113955
$ws.Range("S9").Value = 0.001238678065922237
$ws.Range("T9").Value = 0.001238678065922237
$ws.Range("G10").Value = 0.8410160000000001
$ws.Range("H10").Value = 2.523048
$ws.Range("I10").Value = 0.1349099657534118
$ws.Range("J10").Value = 0.1349099657534118
$ws.Range("M10").Value = 0.5550926666666666
$ws.Range("N10").Value = 1.665278
$ws.Range("O10").Value = 0.1208967663154349
$ws.Range("P10").Value = 0.1208967663154349
$ws.Range("Q10").Value = 0.4668418141493333
$ws.Range("R10").Value = 4.201576327344
$ws.Range("S10").Value = 0.01631017860331355
$ws.Range("T10").Value = 0.01631017860331355
$ws.Range("G11").Value = 0.8410160000000001
$ws.Range("H11").Value = 2.523048
$ws.Range("I11").Value = 0.1349099657534118
$ws.Range("J11").Value = 0.1349099657534118
$ws.Range("O11").Value = 0.7377399926530269
$ws.Range("P11").Value = 0.7377399926530268
$ws.Range("Q11").Value = 2.848776580525334
$ws.Range("R11").Value = 25.638989224728
$ws.Range("S11").Value = 0.09952847714374215
$ws.Range("T11").Value = 0.09952847714374213
$ws.Range("G12").Value = 0.8410160000000001
$ws.Range("H12").Value = 2.523048
$ws.Range("I12").Value = 0.1349099657534118
$ws.Range("J12").Value = 0.1349099657534118
$ws.Range("M12").Value = 0.5311786666666667
$ws.Range("N12").Value = 1.593536
$ws.Range("O12").Value = 0.1156884012202364
$ws.Range("P12").Value = 0.1156884012202364
$ws.Range("Q12").Value = 0.4467297575253334
$ws.Range("R12").Value = 4.020567817728001
$ws.Range("S12").Value = 0.01560751824668906
$ws.Range("T12").Value = 0.01560751824668906
$ws.Range("G13").Value = 0.8410160000000001
$ws.Range("H13").Value = 2.523048
$ws.Range("I13").Value = 0.1349099657534118
$ws.Range("J13").Value = 0.1349099657534118
$ws.Range("M13").Value = 0.117885
$ws.Range("N13").Value = 0.353655
$ws.Range("O13").Value = 0.02567483981130185
$ws.Range("P13").Value = 0.02567483981130185
$ws.Range("Q13").Value = 0.09914317116000002
$ws.Range("R13").Value = 0.8922885404400001
$ws.Range("S13").Value = 0.003463791759667067
$ws.Range("T13").Value = 0.003463791759667067
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.07111199999999999
$ws.Range("H14").Value = 0.213336
$ws.Range("I14").Value = 0.01140729484891681
$ws.Range("J14").Value = 0.01140729484891681
$ws.Range("M14").Value = 0.5550926666666666
$ws.Range("N14").Value = 1.665278
$ws.Range("O14").Value = 0.1208967663154349
$ws.Range("P14").Value = 0.1208967663154349
$ws.Range("Q14").Value = 0.03947374971199999
$ws.Range("R14").Value = 0.3552637474079999
$ws.Range("S14").Value = 0.00137910505964076
$ws.Range("T14").Value = 0.00137910505964076
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.07111199999999999
$ws.Range("H15").Value = 0.213336
$ws.Range("I15").Value = 0.01140729484891681
$ws.Range("J15").Value = 0.01140729484891681
$ws.Range("O15").Value = 0.7377399926530269
$ws.Range("P15").Value = 0.7377399926530268
$ws.Range("Q15").Value = 0.240877938344
$ws.Range("R15").Value = 2.167901445096
$ws.Range("S15").Value = 0.008415617618030799
$ws.Range("T15").Value = 0.008415617618030799
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.07111199999999999
$ws.Range("H16").Value = 0.213336
$ws.Range("I16").Value = 0.01140729484891681
$ws.Range("J16").Value = 0.01140729484891681
$ws.Range("M16").Value = 0.5311786666666667
$ws.Range("N16").Value = 1.593536
$ws.Range("O16").Value = 0.1156884012202364
$ws.Range("P16").Value = 0.1156884012202364
$ws.Range("Q16").Value = 0.037773177344
$ws.Range("R16").Value = 0.339958596096
$ws.Range("S16").Value = 0.001319691703319024
$ws.Range("T16").Value = 0.001319691703319024
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.07111199999999999
$ws.Range("H17").Value = 0.213336
$ws.Range("I17").Value = 0.01140729484891681
$ws.Range("J17").Value = 0.01140729484891681
$ws.Range("M17").Value = 0.117885
$ws.Range("N17").Value = 0.353655
$ws.Range("O17").Value = 0.02567483981130185
$ws.Range("P17").Value = 0.02567483981130185
$ws.Range("Q17").Value = 0.00838303812
$ws.Range("R17").Value = 0.07544734307999999
$ws.Range("S17").Value = 0.0002928804679262279
$ws.Range("T17").Value = 0.0002928804679262279
